$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row in the sheet (data starts at row 2)
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (Changed) date, stored as serial date 45186 (2023-09-17).
# Update every populated row (2..lastRow) to the new serial date 45188 (2023-09-19),
# matching the change recorded in the diff.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
